$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")

# --- Project display-week counter: scrolls the Gantt view back to week 1 ---
$ws.Range("H4").Value = 1

# --- Row 16: category header renamed ---
$ws.Range("B16").Value = "Movements"

# --- Row 17: Motor control task (Lead/Predecessor filled in before the task name) ---
$ws.Range("C17").Value = "M. Amine Gaizi"
$ws.Range("D17").Value = "M. Jeannin"
$ws.Range("B17").Value = "Motor control"
$ws.Range("E17").Value = 43329
$ws.Range("F17").Formula = "=IF(ISBLANK(E17),"" - "",IF(G17=0,E17,E17+G17-1))"
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 15

# --- Row 18: Soldered circuit task ---
$ws.Range("B18").Value = "Soldered circuit"
$ws.Range("C18").Value = "M. Amine Gaizi"
$ws.Range("E18").Value = 43354
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 1

# --- Row 19: Proportional correction task ---
$ws.Range("B19").Value = "Proportional correction"
$ws.Range("C19").Value = "M.Amine Gaizi"
$ws.Range("E19").Value = 43374
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 0.9

# --- Row 20: Wiki Page task ---
$ws.Range("B20").Value = "Wiki Page"
$ws.Range("C20").Value = "M. Amine Gaizi"
$ws.Range("E20").Value = 43364
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 0.6

# --- Cosmetic tweaks matching the author's manual edits ---
$ws.Rows("7").RowHeight = 24.75
$ws.Columns("B").ColumnWidth = 25.6

# --- Final selection left on the sheet ---
[void]$ws.Range("C20").Select()
